# Fill in the "VM Link" column (D) on Sheet1 with the per-participant
# shareable bastion URLs. Rows without a participant name (7, 12, 17, 22)
# are intentionally left blank, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ordered (row, url) pairs -- rows 7, 12, 17, 22 have no participant name
# and are skipped so their VM Link cell stays blank.
$links = @(
    ,(3,  "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/bf6cf0a7-1bc0-4dcb-a0c5-8662a085ea7f")
    ,(4,  "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/374c38c6-88d1-4af1-8d11-b2f66b087746")
    ,(5,  "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/e979dfe7-6554-448c-bbb4-e463b70ec403")
    ,(6,  "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/1a5b1773-0e9e-4831-8bf8-e92587a7bebd")
    ,(8,  "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/817118ff-b149-4bbf-87c8-f3ceae16e0bc")
    ,(9,  "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/b976d368-081d-414e-8432-3fed2b6dfa6c")
    ,(10, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/38687261-e682-4677-8f6d-a0f7ad7633ff")
    ,(11, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/43f93d88-fe8e-4fd7-a6f8-663d62f85696")
    ,(13, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/80ee5d1c-f0c0-4980-a80e-f067619ea8e9")
    ,(14, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/9235c4f6-44ff-44bd-91f3-76bae40279ad")
    ,(15, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/fab8817e-23f4-4e37-be0b-f4f2b56718ba")
    ,(16, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/0b1575b8-3406-4ad7-b5d1-728855490adf")
    ,(18, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/a738499f-5546-4e07-9b1b-eac3af26879d")
    ,(19, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/81f005cb-d609-4ec0-ba23-2cb329e0eef5")
    ,(20, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/46d3df84-28c0-4ff0-98d0-58b192721865")
    ,(21, "https://bst-3662f855-5603-4e55-9e3c-2deb0490560a.bastion.azure.com/api/shareable-url/f06475b3-00be-4797-90f4-1341b7fe8d24")
)

foreach ($pair in $links) {
    $row = $pair[0]
    $url = $pair[1]
    $ws.Cells.Item($row, 4).Value = $url
}

# Leave the selection on the last edited cell, mirroring the saved state.
[void]$ws.Range("D21").Select()
